$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

# The four "level" columns (L1_agg_fuel, L2_CEDS_fuel, L3_agg_sector,
# L4_CEDS_sector) are no longer needed in the Trend_instructions sheet.
# Deleting columns F:I shifts override_normalization/use_as_trend/match_year
# (and their data) left into F:H, and the now-unused shared strings
# (including the stray "x" value in G2) are dropped automatically.
$ws.Range("F1:I1").EntireColumn.Delete()

# Restore the selection that Excel leaves on this sheet after the edit.
$ws.Activate()
$ws.Range("J15").Select()
